$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match header style used by existing header cells (e.g. A1)
$headerRange = $ws.Range("AC1:AE1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows 2-41: Wins, Losses, Ties
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 29).Value = 82
    $ws.Cells.Item($r, 30).Value = 80
    $ws.Cells.Item($r, 31).Value = 0
}
